$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 1104
$ws.Range("I6").Value = 1104
$ws.Range("K6").Value = 3312
$ws.Range("M6").Value = -3200
$ws.Range("H8").Value = 27288.055
$ws.Range("I8").Value = 142956.14
$ws.Range("K8").Value = 428868.42
$ws.Range("M8").Value = -428729.42
$ws.Range("H11").Value = 49.375
$ws.Range("I11").Value = 49.375
$ws.Range("K11").Value = 49.375
$ws.Range("M11").Value = 90.625
$ws.Range("H17").Value = 1997.7812
$ws.Range("J17").Value = 2023.5161
$ws.Range("L17").Value = 6070.5483
$ws.Range("N17").Value = -6406.5483
$ws.Range("H62").Value = 17860208
$ws.Range("I62").Value = 17860208
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 17860208
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -17859584
$ws.Range("H65").Value = 17860208
$ws.Range("I65").Value = 17860208
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 89301040
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -89297920
$ws.Range("H70").Value = 92030.82000000001
$ws.Range("I70").Value = 891.5
$ws.Range("J70").Value = 144110.42
$ws.Range("K70").Value = 2674.5
$ws.Range("L70").Value = 432331.26
$ws.Range("M70").Value = -2404.5
$ws.Range("N70").Value = -432871.26
$ws.Range("H73").Value = 92030.82000000001
$ws.Range("I73").Value = 891.5
$ws.Range("J73").Value = 144110.42
$ws.Range("K73").Value = 2674.5
$ws.Range("L73").Value = 432331.26
$ws.Range("M73").Value = -1738.5
$ws.Range("N73").Value = -434203.26
$ws.Range("H112").Value = 2392.4583
$ws.Range("I112").Value = 1340
$ws.Range("J112").Value = 2438.2173
$ws.Range("K112").Value = 4020
$ws.Range("L112").Value = 7314.651899999999
$ws.Range("M112").Value = -2912
$ws.Range("N112").Value = -9530.651899999999
$ws.Range("H132").Value = 1818.4286
$ws.Range("I132").Value = 1052.575
$ws.Range("K132").Value = 3157.725
$ws.Range("M132").Value = -627.7250000000004
$ws.Range("H133").Value = 58076.734
$ws.Range("J133").Value = 58076.734
$ws.Range("L133").Value = 58076.734
$ws.Range("N133").Value = -68196.734
$ws.Range("N62").Value = $null
$ws.Range("N65").Value = $null

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2530.7317
$ws.Range("I32").Value = 2175.5066
$ws.Range("J32").Value = 8001.2
$ws.Range("K32").Value = 2175.5066
$ws.Range("L32").Value = 8001.2
$ws.Range("M32").Value = -1888.5066
$ws.Range("N32").Value = -8575.200000000001
$ws.Range("H45").Value = 3687.9443
$ws.Range("I45").Value = 2884.7144
$ws.Range("J45").Value = 6499.25
$ws.Range("K45").Value = 2884.7144
$ws.Range("L45").Value = 6499.25
$ws.Range("M45").Value = -2507.7144
$ws.Range("N45").Value = -7253.25
$ws.Range("H74").Value = 1710.6666
$ws.Range("I74").Value = 1385.1666
$ws.Range("K74").Value = 1385.1666
$ws.Range("M74").Value = -511.1666
$ws.Range("H77").Value = 1710.6666
$ws.Range("I77").Value = 1385.1666
$ws.Range("K77").Value = 6925.833000000001
$ws.Range("M77").Value = -2557.833000000001
$ws.Range("H110").Value = 418446.66
$ws.Range("I110").Value = 501636
$ws.Range("K110").Value = 501636
$ws.Range("M110").Value = -499591
$ws.Range("H124").Value = 55769.145
$ws.Range("J124").Value = 55769.145
$ws.Range("L124").Value = 55769.145
$ws.Range("N124").Value = -65589.14499999999
$ws.Range("H125").Value = 54944.5
$ws.Range("J125").Value = 54944.5
$ws.Range("L125").Value = 54944.5
$ws.Range("N125").Value = -64784.5
$ws.Range("H132").Value = 3333.3027
$ws.Range("I132").Value = 1125.4529
$ws.Range("K132").Value = 3376.3587
$ws.Range("M132").Value = -846.3586999999998

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 250
$ws.Range("I5").Value = 100
$ws.Range("K5").Value = 100
$ws.Range("M5").Value = 13
$ws.Range("H134").Value = 3731.28
$ws.Range("I134").Value = 2488.9707
$ws.Range("J134").Value = 6371.1875
$ws.Range("K134").Value = 7466.9121
$ws.Range("L134").Value = 19113.5625
$ws.Range("M134").Value = -4931.9121
$ws.Range("N134").Value = -24183.5625

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 80000000
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("H16").Value = 3056
$ws.Range("J16").Value = 3999.6667
$ws.Range("L16").Value = 3999.6667
$ws.Range("N16").Value = -4573.6667
$ws.Range("H31").Value = 347765.06
$ws.Range("I31").Value = 527835.2
$ws.Range("K31").Value = 527835.2
$ws.Range("M31").Value = -527540.2
$ws.Range("H34").Value = 347765.06
$ws.Range("I34").Value = 527835.2
$ws.Range("K34").Value = 527835.2
$ws.Range("M34").Value = -527633.2
$ws.Range("H58").Value = 199036.22
$ws.Range("I58").Value = 313835.6
$ws.Range("J58").Value = 5689.8945
$ws.Range("K58").Value = 313835.6
$ws.Range("L58").Value = 5689.8945
$ws.Range("M58").Value = -313632.6
$ws.Range("N58").Value = -6095.8945
$ws.Range("H62").Value = 3624.5
$ws.Range("I62").Value = 3624.5
$ws.Range("K62").Value = 3624.5
$ws.Range("M62").Value = -3000.5
$ws.Range("H65").Value = 3624.5
$ws.Range("I65").Value = 3624.5
$ws.Range("K65").Value = 18122.5
$ws.Range("M65").Value = -15002.5
$ws.Range("H113").Value = 3056
$ws.Range("J113").Value = 3999.6667
$ws.Range("L113").Value = 3999.6667
$ws.Range("N113").Value = -8339.6667
$ws.Range("H122").Value = 3048.5715
$ws.Range("I122").Value = 2353.7334
$ws.Range("J122").Value = 4785.6665
$ws.Range("K122").Value = 7061.2002
$ws.Range("L122").Value = 14356.9995
$ws.Range("M122").Value = -4611.2002
$ws.Range("N122").Value = -19256.9995
$ws.Range("H136").Value = 199036.22
$ws.Range("I136").Value = 313835.6
$ws.Range("J136").Value = 5689.8945
$ws.Range("K136").Value = 941506.7999999999
$ws.Range("L136").Value = 17069.6835
$ws.Range("M136").Value = -938956.7999999999
$ws.Range("N136").Value = -22169.6835
$ws.Range("M4").Value = $null

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 80946.3
$ws.Range("J5").Value = 1477.1666
$ws.Range("L5").Value = 4431.4998
$ws.Range("N5").Value = -4655.4998
$ws.Range("H39").Value = 7016.2144
$ws.Range("I39").Value = 939.6
$ws.Range("J39").Value = 10392.111
$ws.Range("K39").Value = 2818.8
$ws.Range("L39").Value = 31176.333
$ws.Range("M39").Value = -2524.8
$ws.Range("N39").Value = -31764.333
$ws.Range("H127").Value = 1944
$ws.Range("J127").Value = 1944
$ws.Range("L127").Value = 5832
$ws.Range("N127").Value = -15752
$ws.Range("H131").Value = 2985.8044
$ws.Range("I131").Value = 1115.6666
$ws.Range("K131").Value = 3346.9998
$ws.Range("M131").Value = 1693.0002
$ws.Range("H135").Value = 80946.3
$ws.Range("J135").Value = 1477.1666
$ws.Range("L135").Value = 13294.4994
$ws.Range("N135").Value = -18364.4994

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 1880375.2
$ws.Range("I3").Value = 7500499.5
$ws.Range("J3").Value = 7000.5
$ws.Range("K3").Value = 7500499.5
$ws.Range("L3").Value = 7000.5
$ws.Range("M3").Value = -7500383.5
$ws.Range("N3").Value = -7232.5
$ws.Range("H138").Value = 49800
$ws.Range("J138").Value = 49800
$ws.Range("L138").Value = 49800
$ws.Range("N138").Value = -60080

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 5167.5293
$ws.Range("I46").Value = 4622.727
$ws.Range("J46").Value = 6166.3335
$ws.Range("K46").Value = 4622.727
$ws.Range("L46").Value = 6166.3335
$ws.Range("M46").Value = -4434.727
$ws.Range("N46").Value = -6542.3335
$ws.Range("H61").Value = 3552.6562
$ws.Range("I61").Value = 2767.28
$ws.Range("K61").Value = 2767.28
$ws.Range("M61").Value = -2565.28
$ws.Range("H68").Value = 394
$ws.Range("I68").Value = 394
$ws.Range("K68").Value = 394
$ws.Range("M68").Value = 355
$ws.Range("H71").Value = 394
$ws.Range("I71").Value = 394
$ws.Range("K71").Value = 1970
$ws.Range("M71").Value = 1774
$ws.Range("H82").Value = 2136.889
$ws.Range("I82").Value = 1720.125
$ws.Range("J82").Value = 2470.3
$ws.Range("K82").Value = 1720.125
$ws.Range("L82").Value = 2470.3
$ws.Range("M82").Value = -1359.125
$ws.Range("N82").Value = -3192.3
$ws.Range("H85").Value = 2136.889
$ws.Range("I85").Value = 1720.125
$ws.Range("J85").Value = 2470.3
$ws.Range("K85").Value = 1720.125
$ws.Range("L85").Value = 2470.3
$ws.Range("M85").Value = -472.125
$ws.Range("N85").Value = -4966.3
$ws.Range("H113").Value = 3552.6562
$ws.Range("I113").Value = 2767.28
$ws.Range("K113").Value = 2767.28
$ws.Range("M113").Value = -597.2800000000002
$ws.Range("H127").Value = 75000
$ws.Range("J127").Value = 75000
$ws.Range("L127").Value = 75000
$ws.Range("N127").Value = -84920
$ws.Range("H132").Value = 5072.9375
$ws.Range("I132").Value = 4240.8887
$ws.Range("K132").Value = 12722.6661
$ws.Range("M132").Value = -10192.6661

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 7322.2856
$ws.Range("I62").Value = 8250
$ws.Range("K62").Value = 8250
$ws.Range("M62").Value = -7626
$ws.Range("H65").Value = 7322.2856
$ws.Range("I65").Value = 8250
$ws.Range("K65").Value = 41250
$ws.Range("M65").Value = -38130
$ws.Range("H113").Value = 1470.0416
$ws.Range("I113").Value = 1127.5555
$ws.Range("J113").Value = 2497.5
$ws.Range("K113").Value = 3382.6665
$ws.Range("L113").Value = 7492.5
$ws.Range("M113").Value = -1212.6665
$ws.Range("N113").Value = -11832.5
